$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the U/V totals formulas on rows 8 and 9 ---
# U8/U9: use S (sous-total "Heures Sup") instead of the stray P reference
$ws.Range("U8").Formula = "=I8+N8+S8"
$ws.Range("U9").Formula = "=I9+N9+S9"
# V8/V9: use T instead of S (since U now already uses S)
$ws.Range("V8").Formula = "=J8+O8+T8"
$ws.Range("V9").Formula = "=J9+O9+T9"

# --- Row 11: add the missing bordered, empty S11 cell (same style as neighbours) ---
$ws.Range("S11").Borders.LineStyle = 1

# --- Row 11: these cells held a stray "=  " formula evaluating to 0; clear them to plain empty cells ---
$ws.Range("J11").ClearContents()
$ws.Range("O11").ClearContents()
$ws.Range("T11").ClearContents()
$ws.Range("U11").ClearContents()
$ws.Range("V11").ClearContents()

# --- Row 12: add the missing bordered, empty S12 cell (same style as neighbours) ---
$ws.Range("S12").Borders.LineStyle = 1
